$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 update
$ws.Range("B2").Value = 73834

# Row 3 and Row 4 "swap" species data (D,E,F,G,H,Q,R) and the record id (A),
# while B gets new (non-swapped) values.
$ws.Range("A3").Value = 112491430
$ws.Range("B3").Value = 90814
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4364
$ws.Range("F3").Value = "Dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum ferrugineum"
$ws.Range("H3").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q3").Value = 356556
$ws.Range("R3").Value = 6742347

$ws.Range("A4").Value = 112491413
$ws.Range("B4").Value = 90830
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 2059
$ws.Range("F4").Value = "Skrovlig taggsvamp"
$ws.Range("G4").Value = "Hydnellum scabrosum"
$ws.Range("H4").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q4").Value = 356670
$ws.Range("R4").Value = 6742658

# Row 5: B5 update
$ws.Range("B5").Value = 90857
